$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 14 -> P13
$ws.Range("B14").Value = 26
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0

# Row 15 -> P14
$ws.Range("B15").Value = 34
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0

$ws.Range("I14").Select()
